$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37, shifting existing rows 37-65 down to 38-66.
$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the new weekly record (copy of the row that
# used to be at 37, with updated date / price columns).
$ws.Cells.Item(37, 1).Value = 11
$ws.Cells.Item(37, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value = "Bíobío"
$ws.Cells.Item(37, 4).Value = 44447
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = 100112043
$ws.Cells.Item(37, 7).Value = "Pepino ensalada"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 16000
$ws.Cells.Item(37, 12).Value = 17000
$ws.Cells.Item(37, 13).Value = 16500
$ws.Cells.Item(37, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 275
$ws.Cells.Item(37, 17).Value = 60
$ws.Cells.Item(37, 18).Value = "Hortaliza"
